# The document contains four "<id>...</id>" tag sequences that were each
# split across three runs (open-tag run / value run / close-tag run) because
# the id value was typed/edited separately from the tag markup. The commit
# re-downloads/re-normalizes the tei markup so each "<id>VALUE</id>" becomes
# a single contiguous run (taking on the formatting of the opening "<id>"
# run, i.e. Courier New / color 7f6000).
#
# A plain Find & Replace over the exact (unique) text of each occurrence
# collapses the matched runs into one run using the formatting of the first
# character of the match, which reproduces that exact target shape.

$d = $word.ActiveDocument

$ids = @("p132r_2", "p132v_1", "p132v_2", "p132v_3")

foreach ($val in $ids) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $old = "<id>" + $val + "</id>"
    $new = "<id>" + $val + "</id>"
    $result = $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    Write-Host "Replaced $val -> $result"
}
